# Hortaliza, Vega Modelo de Temuco - Caigua: weekly fruit/vegetable price refresh.
# The reporting date (D) and volume/price figures (J, K, L, M, P) for each
# data row (2-26) are refreshed with the latest values for that market row,
# which in effect redistributes the previously-recorded rows' data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = 44837; J = 80;  K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 3;  D = 44839; J = 80;  K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 4;  D = 44838; J = 10;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 5;  D = 44819; J = 100; K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 6;  D = 44749; J = 50;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 7;  D = 44826; J = 50;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 8;  D = 44508; J = 40;  K = 10000; L = 10000; M = 10000; P = 667 },
    @{ Row = 9;  D = 44827; J = 20;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 10; D = 44845; J = 20;  K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 11; D = 44767; J = 50;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 12; D = 44756; J = 80;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 13; D = 44841; J = 20;  K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 14; D = 44771; J = 40;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 15; D = 44811; J = 30;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 16; D = 44757; J = 30;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 17; D = 44525; J = 40;  K = 8000;  L = 8000;  M = 8000;  P = 533 },
    @{ Row = 18; D = 44518; J = 50;  K = 10000; L = 10000; M = 10000; P = 667 },
    @{ Row = 19; D = 44824; J = 20;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 20; D = 44812; J = 80;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 21; D = 44755; J = 50;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 22; D = 44830; J = 25;  K = 12000; L = 12000; M = 12000; P = 800 },
    @{ Row = 23; D = 44825; J = 30;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 24; D = 44776; J = 80;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 25; D = 44813; J = 20;  K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 26; D = 44769; J = 50;  K = 20000; L = 20000; M = 20000; P = 1333 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value  = $u.D   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $u.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $u.K   # K: Precio mínimo
    $ws.Cells.Item($r, 12).Value = $u.L   # L: Precio máximo
    $ws.Cells.Item($r, 13).Value = $u.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $u.P   # P: Precio $/Kg
}
